$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list with new Price / Volume(1h) figures.
# Values must stay as text (to match the original inlineStr cells),
# so force a Text number format before writing each cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.67%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.04%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.704"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.74%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06209"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.45%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.727"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.47%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8513"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.91%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9077"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.51%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.01%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04707"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-10.77%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07097"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.28%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03175"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.98%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09061"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.79%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001540"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.85%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006180"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.06%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006022"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.04%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.470"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.45%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.170"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.24%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.32%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.51%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.112"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.60%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04220"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.60%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.10%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004115"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.87%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.10%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03899"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.66%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1114"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.28%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004132"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.12%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.70%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01345"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-9.89%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005175"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.30%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.10%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.03504"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-35.75%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1589"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "17.44%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.10%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.10%"
